$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text (country name) changes from re-ranking ---
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 7 de Julio de 2020 a las 18:47'
$ws.Cells.Item(30, 1).Value = 'Irak'
$ws.Cells.Item(31, 1).Value = 'Bielorrusia'
$ws.Cells.Item(32, 1).Value = 'Ecuador'
$ws.Cells.Item(57, 1).Value = 'Azerbaiyan'
$ws.Cells.Item(58, 1).Value = 'Ghana'
$ws.Cells.Item(62, 1).Value = 'Argelia'
$ws.Cells.Item(63, 1).Value = 'Serbia'
$ws.Cells.Item(116, 1).Value = 'Zambia'
$ws.Cells.Item(117, 1).Value = 'Islandia'
$ws.Cells.Item(118, 1).Value = 'Lituania'
$ws.Cells.Item(119, 1).Value = 'Malaui'
$ws.Cells.Item(120, 1).Value = 'Guinea-Bisau'
$ws.Cells.Item(121, 1).Value = 'Eslovaquia'
$ws.Cells.Item(122, 1).Value = 'Eslovenia'
$ws.Cells.Item(129, 1).Value = 'Tunez'
$ws.Cells.Item(130, 1).Value = 'Benin'
$ws.Cells.Item(143, 1).Value = 'Montenegro'
$ws.Cells.Item(144, 1).Value = 'Liberia'
$ws.Cells.Item(145, 1).Value = 'Republica del Chad'
$ws.Cells.Item(146, 1).Value = 'Principado de Andorra'
$ws.Cells.Item(209, 1).Value = 'Islas Malvinas'
$ws.Cells.Item(210, 1).Value = 'Groenlandia'

# --- Updated case numbers ---
# Row 4
$ws.Cells.Item(4, 2).Value = 3057011
$ws.Cells.Item(4, 3).Value = 16178
$ws.Cells.Item(4, 4).Value = 1326770
$ws.Cells.Item(4, 5).Value = 1596939
$ws.Cells.Item(4, 7).Value = 323
$ws.Cells.Item(4, 8).Value = 133302
# Row 5
$ws.Cells.Item(5, 2).Value = 1643539
$ws.Cells.Item(5, 3).Value = 17468
$ws.Cells.Item(5, 5).Value = 505217
$ws.Cells.Item(5, 7).Value = 537
$ws.Cells.Item(5, 8).Value = 66093
# Row 6
$ws.Cells.Item(6, 2).Value = 740131
$ws.Cells.Item(6, 3).Value = 19785
$ws.Cells.Item(6, 4).Value = 455191
$ws.Cells.Item(6, 5).Value = 264304
$ws.Cells.Item(6, 7).Value = 462
$ws.Cells.Item(6, 8).Value = 20636
# Row 9
$ws.Cells.Item(9, 4).Value = 268245
$ws.Cells.Item(9, 5).Value = 26340
# Row 11
$ws.Cells.Item(11, 2).Value = 286349
$ws.Cells.Item(11, 3).Value = 581
$ws.Cells.Item(11, 7).Value = 155
$ws.Cells.Item(11, 8).Value = 44391
# Row 14
$ws.Cells.Item(14, 2).Value = 241956
$ws.Cells.Item(14, 3).Value = 137
$ws.Cells.Item(14, 4).Value = 192815
$ws.Cells.Item(14, 5).Value = 14242
$ws.Cells.Item(14, 7).Value = 30
$ws.Cells.Item(14, 8).Value = 34899
# Row 17
$ws.Cells.Item(17, 2).Value = 207897
$ws.Cells.Item(17, 3).Value = 1053
$ws.Cells.Item(17, 4).Value = 185292
$ws.Cells.Item(17, 5).Value = 17345
$ws.Cells.Item(17, 7).Value = 19
$ws.Cells.Item(17, 8).Value = 5260
# Row 19
$ws.Cells.Item(19, 2).Value = 198172
$ws.Cells.Item(19, 3).Value = 115
$ws.Cells.Item(19, 5).Value = 6379
# Row 30
$ws.Cells.Item(30, 2).Value = 64701
$ws.Cells.Item(30, 3).Value = 2426
$ws.Cells.Item(30, 4).Value = 36252
$ws.Cells.Item(30, 5).Value = 25764
$ws.Cells.Item(30, 7).Value = 118
$ws.Cells.Item(30, 8).Value = 2685
# Row 31
$ws.Cells.Item(31, 2).Value = 64003
$ws.Cells.Item(31, 3).Value = 199
$ws.Cells.Item(31, 4).Value = 51902
$ws.Cells.Item(31, 5).Value = 11665
$ws.Cells.Item(31, 7).Value = 7
$ws.Cells.Item(31, 8).Value = 436
# Row 32
$ws.Cells.Item(32, 2).Value = 62380
$ws.Cells.Item(32, 4).Value = 28872
$ws.Cells.Item(32, 5).Value = 28687
$ws.Cells.Item(32, 8).Value = 4821
# Row 41
$ws.Cells.Item(41, 4).Value = 41002
$ws.Cells.Item(41, 5).Value = 4112
# Row 45
$ws.Cells.Item(45, 2).Value = 38430
$ws.Cells.Item(45, 3).Value = 302
$ws.Cells.Item(45, 4).Value = 19564
$ws.Cells.Item(45, 5).Value = 18045
$ws.Cells.Item(45, 7).Value = 17
$ws.Cells.Item(45, 8).Value = 821
# Row 49
$ws.Cells.Item(49, 2).Value = 31886
$ws.Cells.Item(49, 3).Value = 1137
$ws.Cells.Item(49, 4).Value = 18192
$ws.Cells.Item(49, 5).Value = 13352
$ws.Cells.Item(49, 7).Value = 8
$ws.Cells.Item(49, 8).Value = 342
# Row 57
$ws.Cells.Item(57, 2).Value = 21374
$ws.Cells.Item(57, 3).Value = 537
$ws.Cells.Item(57, 4).Value = 12635
$ws.Cells.Item(57, 5).Value = 8474
$ws.Cells.Item(57, 7).Value = 7
$ws.Cells.Item(57, 8).Value = 265
# Row 58
$ws.Cells.Item(58, 2).Value = 21077
$ws.Cells.Item(58, 4).Value = 16070
$ws.Cells.Item(58, 5).Value = 4878
$ws.Cells.Item(58, 8).Value = 129
# Row 61
$ws.Cells.Item(61, 2).Value = 18141
$ws.Cells.Item(61, 3).Value = 235
$ws.Cells.Item(61, 5).Value = 6297
# Row 62
$ws.Cells.Item(62, 2).Value = 16879
$ws.Cells.Item(62, 3).Value = 475
$ws.Cells.Item(62, 4).Value = 11884
$ws.Cells.Item(62, 5).Value = 4027
$ws.Cells.Item(62, 7).Value = 9
$ws.Cells.Item(62, 8).Value = 968
# Row 63
$ws.Cells.Item(63, 2).Value = 16719
$ws.Cells.Item(63, 3).Value = 299
$ws.Cells.Item(63, 4).Value = 13366
$ws.Cells.Item(63, 5).Value = 3023
$ws.Cells.Item(63, 7).Value = 13
$ws.Cells.Item(63, 8).Value = 330
# Row 69
$ws.Cells.Item(69, 2).Value = 12639
$ws.Cells.Item(69, 3).Value = 73
$ws.Cells.Item(69, 5).Value = 4415
$ws.Cells.Item(69, 7).Value = 1
$ws.Cells.Item(69, 8).Value = 351
# Row 72
$ws.Cells.Item(72, 2).Value = 9997
$ws.Cells.Item(72, 3).Value = 103
$ws.Cells.Item(72, 4).Value = 5034
$ws.Cells.Item(72, 5).Value = 4341
$ws.Cells.Item(72, 7).Value = 6
$ws.Cells.Item(72, 8).Value = 622
# Row 73
$ws.Cells.Item(73, 2).Value = 8941
$ws.Cells.Item(73, 3).Value = 5
$ws.Cells.Item(73, 5).Value = 552
# Row 96
$ws.Cells.Item(96, 2).Value = 4603
$ws.Cells.Item(96, 3).Value = 61
$ws.Cells.Item(96, 4).Value = 4056
$ws.Cells.Item(96, 5).Value = 437
# Row 99
$ws.Cells.Item(99, 2).Value = 3589
$ws.Cells.Item(99, 3).Value = 27
$ws.Cells.Item(99, 5).Value = 2022
$ws.Cells.Item(99, 7).Value = 1
$ws.Cells.Item(99, 8).Value = 193
# Row 112
$ws.Cells.Item(112, 2).Value = 2081
$ws.Cells.Item(112, 3).Value = 4
$ws.Cells.Item(112, 5).Value = 115
# Row 115
$ws.Cells.Item(115, 2).Value = 1907
$ws.Cells.Item(115, 3).Value = 22
$ws.Cells.Item(115, 4).Value = 1348
$ws.Cells.Item(115, 5).Value = 523
# Row 116
$ws.Cells.Item(116, 2).Value = 1895
$ws.Cells.Item(116, 3).Value = 263
$ws.Cells.Item(116, 4).Value = 1348
$ws.Cells.Item(116, 5).Value = 505
$ws.Cells.Item(116, 7).Value = 12
$ws.Cells.Item(116, 8).Value = 42
# Row 117
$ws.Cells.Item(117, 2).Value = 1873
$ws.Cells.Item(117, 3).Value = 7
$ws.Cells.Item(117, 4).Value = 1847
$ws.Cells.Item(117, 5).Value = 16
$ws.Cells.Item(117, 8).Value = 10
# Row 118
$ws.Cells.Item(118, 2).Value = 1844
$ws.Cells.Item(118, 3).Value = 3
$ws.Cells.Item(118, 4).Value = 1547
$ws.Cells.Item(118, 5).Value = 218
$ws.Cells.Item(118, 8).Value = 79
# Row 119
$ws.Cells.Item(119, 2).Value = 1818
$ws.Cells.Item(119, 3).Value = 76
$ws.Cells.Item(119, 4).Value = 317
$ws.Cells.Item(119, 5).Value = 1482
$ws.Cells.Item(119, 8).Value = 19
# Row 120
$ws.Cells.Item(120, 2).Value = 1790
$ws.Cells.Item(120, 3).Value = 0
$ws.Cells.Item(120, 4).Value = 760
$ws.Cells.Item(120, 5).Value = 1005
$ws.Cells.Item(120, 8).Value = 25
# Row 121
$ws.Cells.Item(121, 2).Value = 1767
$ws.Cells.Item(121, 3).Value = 2
$ws.Cells.Item(121, 4).Value = 1473
$ws.Cells.Item(121, 5).Value = 266
$ws.Cells.Item(121, 8).Value = 28
# Row 122
$ws.Cells.Item(122, 2).Value = 1739
$ws.Cells.Item(122, 3).Value = 23
$ws.Cells.Item(122, 4).Value = 1423
$ws.Cells.Item(122, 5).Value = 205
$ws.Cells.Item(122, 8).Value = 111
# Row 129
$ws.Cells.Item(129, 2).Value = 1205
$ws.Cells.Item(129, 3).Value = 6
$ws.Cells.Item(129, 4).Value = 1049
$ws.Cells.Item(129, 5).Value = 106
$ws.Cells.Item(129, 8).Value = 50
# Row 130
$ws.Cells.Item(130, 4).Value = 333
$ws.Cells.Item(130, 5).Value = 845
$ws.Cells.Item(130, 8).Value = 21
# Row 143
$ws.Cells.Item(143, 2).Value = 907
$ws.Cells.Item(143, 3).Value = 66
$ws.Cells.Item(143, 4).Value = 315
$ws.Cells.Item(143, 5).Value = 575
$ws.Cells.Item(143, 7).Value = 3
$ws.Cells.Item(143, 8).Value = 17
# Row 144
$ws.Cells.Item(144, 2).Value = 891
$ws.Cells.Item(144, 4).Value = 377
$ws.Cells.Item(144, 5).Value = 475
$ws.Cells.Item(144, 8).Value = 39
# Row 145
$ws.Cells.Item(145, 2).Value = 872
$ws.Cells.Item(145, 4).Value = 787
$ws.Cells.Item(145, 5).Value = 11
$ws.Cells.Item(145, 8).Value = 74
# Row 146
$ws.Cells.Item(146, 2).Value = 855
$ws.Cells.Item(146, 4).Value = 800
$ws.Cells.Item(146, 5).Value = 3
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 52
# Row 155
$ws.Cells.Item(155, 5).Value = 75
$ws.Cells.Item(155, 7).Value = 1
$ws.Cells.Item(155, 8).Value = 3
